$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing C4 value (2 -> 3)
$ws.Range("C4").Value = 3

# Add new row 5 - copy date formatting from A4, then set values
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 44902
$ws.Range("B5").Value = "Starting on the spatial ACF"
$ws.Range("C5").Value = 1.5

# Add new row 6
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 44902
$ws.Range("B6").Value = "Mapping the overall data trends"
$ws.Range("C6").Value = 1.5

$excel.CutCopyMode = $false

# Update the selection to match the diff (C7)
$ws.Range("C7").Select()
